# Update "想去人数" (interest count) figures in column F across the
# "展览" (Exhibitions), "演出" (Performances) and "全部类型" (All types)
# sheets, reflecting freshly scraped counts (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Range("F3").Value  = 247
$wsExpo.Range("F4").Value  = 266
$wsExpo.Range("F8").Value  = 2222
$wsExpo.Range("F9").Value  = 368
$wsExpo.Range("F12").Value = 83
$wsExpo.Range("F13").Value = 2550
$wsExpo.Range("F15").Value = 1350
$wsExpo.Range("F16").Value = 4706
$wsExpo.Range("F18").Value = 5135
$wsExpo.Range("F19").Value = 1685
$wsExpo.Range("F20").Value = 2877
$wsExpo.Range("F21").Value = 3276
$wsExpo.Range("F23").Value = 1558
$wsExpo.Range("F24").Value = 259
$wsExpo.Range("F28").Value = 1009
$wsExpo.Range("F29").Value = 1884
$wsExpo.Range("F31").Value = 282
$wsExpo.Range("F32").Value = 718
$wsExpo.Range("F33").Value = 159
$wsExpo.Range("F34").Value = 338
$wsExpo.Range("F35").Value = 418

# --- Sheet 2: 演出 (Performances) ---
$wsShow = $wb.Worksheets.Item(2)
$wsShow.Range("F3").Value  = 99
$wsShow.Range("F11").Value = 195

# --- Sheet 4: 全部类型 (All types, merges the above two sheets) ---
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F3").Value  = 99
$wsAll.Range("F8").Value  = 247
$wsAll.Range("F10").Value = 266
$wsAll.Range("F13").Value = 2222
$wsAll.Range("F14").Value = 368
$wsAll.Range("F19").Value = 83
$wsAll.Range("F21").Value = 2550
$wsAll.Range("F22").Value = 1350
$wsAll.Range("F23").Value = 195
$wsAll.Range("F26").Value = 4706
$wsAll.Range("F28").Value = 5135
$wsAll.Range("F29").Value = 1685
$wsAll.Range("F30").Value = 2877
$wsAll.Range("F31").Value = 3276
$wsAll.Range("F35").Value = 1558
$wsAll.Range("F37").Value = 259
$wsAll.Range("F41").Value = 1009
$wsAll.Range("F43").Value = 1884
$wsAll.Range("F45").Value = 282
$wsAll.Range("F46").Value = 718
$wsAll.Range("F47").Value = 159
$wsAll.Range("F48").Value = 338
$wsAll.Range("F49").Value = 418
